$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure text columns keep their exact string representation (avoid Excel auto-numeric conversion)
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '37.966.23'
$ws.Range("E2").Value = '  -0.08%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.035.61'
$ws.Range("E3").Value = '  -0.95%  '
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '228.50'
$ws.Range("E5").Value = '  -0.70%  '
$ws.Range("E6").Value = '  -0.85%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '60.47'
$ws.Range("E7").Value = '  +2.78%  '
$ws.Range("E8").Value = '  -0.03%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.380'
$ws.Range("E9").Value = '  -1.64%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0822'
$ws.Range("E10").Value = '  +0.89%  '
$ws.Range("E11").Value = '  +0.33%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '2.337.75'
$ws.Range("E12").Value = '  -0.91%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '14.54'
$ws.Range("E13").Value = '  -0.68%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '21.31'
$ws.Range("E14").Value = '  +1.68%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '0.763'
$ws.Range("E15").Value = '  +1.15%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '5.17'
$ws.Range("E16").Value = '  -2.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.041.47'
$ws.Range("E17").Value = '  -0.66%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '37.896.96'
$ws.Range("E18").Value = '  -0.03%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '69.72'
$ws.Range("E19").Value = '  -0.08%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '5.89'
$ws.Range("E20").Value = '  -7.38%  '
$ws.Range("E21").Value = '  -1.49%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '224.23'
$ws.Range("E22").Value = '  -0.20%  '
$ws.Range("E23").Value = '  +0.07%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '2.42'
$ws.Range("E24").Value = '  -0.39%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.25'
$ws.Range("E25").Value = '  -0.44%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '9.31'
$ws.Range("E26").Value = '  +0.11%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '166.67'
$ws.Range("E27").Value = '  +0.04%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.131'
$ws.Range("E28").Value = '  -0.90%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '18.86'
$ws.Range("E29").Value = '  -0.96%  '
$ws.Range("E30").Value = '  -3.72%  '
$ws.Range("E31").Value = '  +0.56%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '2.24'
$ws.Range("E32").Value = '  +8.32%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '4.40'
$ws.Range("E33").Value = '  -2.94%  '
$ws.Range("B34").Value = 'Hedera'
$ws.Range("C34").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '0.0606'
$ws.Range("E34").Value = '  -0.44%  '
$ws.Range("B35").Value = 'InternetComputer(DFINITY)'
$ws.Range("C35").Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.52'
$ws.Range("E35").Value = '  -1.93%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.44'
$ws.Range("E36").Value = '  +6.20%  '
$ws.Range("E37").Value = '  -2.46%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.28'
$ws.Range("E38").Value = '  +0.05%  '
$ws.Range("E39").Value = '  +0.08%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '17.70'
$ws.Range("E40").Value = '  +6.96%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.535.43'
$ws.Range("E41").Value = '  +1.89%  '
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '96.22'
$ws.Range("E43").Value = '  -0.87%  '
$ws.Range("E44").Value = '  -2.17%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0913'
$ws.Range("E45").Value = '  -0.78%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '1.10'
$ws.Range("E46").Value = '  -2.50%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.01'
$ws.Range("E47").Value = '  -2.40%  '
$ws.Range("E48").Value = '  -0.82%  '
$ws.Range("E49").Value = '  +0.04%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '7.09'
$ws.Range("E50").Value = '  -0.16%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.227.29'
$ws.Range("E51").Value = '  -0.81%  '
